# Update "想去人数" (F column) counts across all sheets to reflect the
# latest generated output (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (row, newValue) pairs for column F.
$updates = @{
    "展览" = @(
        @{Row=4;  Value=317},
        @{Row=5;  Value=5702},
        @{Row=7;  Value=7703},
        @{Row=10; Value=3873},
        @{Row=11; Value=76},
        @{Row=14; Value=206},
        @{Row=18; Value=108},
        @{Row=20; Value=617},
        @{Row=21; Value=3904},
        @{Row=24; Value=5337},
        @{Row=25; Value=441},
        @{Row=26; Value=2116},
        @{Row=28; Value=356},
        @{Row=29; Value=7949},
        @{Row=30; Value=33},
        @{Row=33; Value=2203},
        @{Row=35; Value=1305},
        @{Row=36; Value=21},
        @{Row=41; Value=1181},
        @{Row=42; Value=1177},
        @{Row=44; Value=1337},
        @{Row=45; Value=2098},
        @{Row=46; Value=132},
        @{Row=47; Value=226},
        @{Row=48; Value=1219}
    )
    "演出" = @(
        @{Row=11; Value=124},
        @{Row=15; Value=12}
    )
    "本地生活" = @(
        @{Row=2; Value=579},
        @{Row=4; Value=67}
    )
    "全部类型" = @(
        @{Row=4;  Value=67},
        @{Row=5;  Value=317},
        @{Row=6;  Value=5702},
        @{Row=8;  Value=3873},
        @{Row=9;  Value=76},
        @{Row=15; Value=108},
        @{Row=18; Value=617},
        @{Row=19; Value=3904},
        @{Row=23; Value=5337},
        @{Row=24; Value=441},
        @{Row=25; Value=2116},
        @{Row=27; Value=356},
        @{Row=28; Value=7949},
        @{Row=29; Value=33},
        @{Row=31; Value=2203},
        @{Row=33; Value=1305},
        @{Row=38; Value=1181},
        @{Row=39; Value=1177},
        @{Row=42; Value=1337},
        @{Row=44; Value=2098},
        @{Row=45; Value=132},
        @{Row=46; Value=226},
        @{Row=49; Value=1219}
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Range("F" + $entry.Row).Value = $entry.Value
    }
}
